$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CottonObserved")

# --- New header cells (row 1), columns AG:AN ---------------------------------
# Written in an order that reproduces the target sharedStrings insertion
# order: Seed.N, Seed.Nerror, Seed.NConc, Seed.NConcError, Burr.N, Burr.NError,
# Burr.NConcError, Burr.NConc  (note AN is written before AM).
$ws.Range("AG1").Value = "Cotton.Seed.N"
$ws.Range("AH1").Value = "Cotton.Seed.Nerror"
$ws.Range("AI1").Value = "Cotton.Seed.NConc"
$ws.Range("AJ1").Value = "Cotton.Seed.NConcError"
$ws.Range("AK1").Value = "Cotton.Burr.N"
$ws.Range("AL1").Value = "Cotton.Burr.NError"
$ws.Range("AN1").Value = "Cotton.Burr.NConcError"
$ws.Range("AM1").Value = "Cotton.Burr.NConc"

# --- New data cells ------------------------------------------------------
# Row 21
$ws.Range("AG21").Value = [double]"7.6506823406061057"
$ws.Range("AH21").Value = [double]"0.87325237998046579"
$ws.Range("AI21").Value = [double]"4.0004999999999999E-2"
$ws.Range("AJ21").Value = [double]"1.88860618799512E-3"
$ws.Range("AK21").Value = [double]"2.5361029191403284"
$ws.Range("AL21").Value = [double]"0.31169440130238801"
$ws.Range("AM21").Value = [double]"1.2669999999999999E-2"
$ws.Range("AN21").Value = [double]"2.12967916206487E-3"

# Row 23
$ws.Range("AI23").Value = [double]"4.3422500000000003E-2"
$ws.Range("AJ23").Value = [double]"2.2055290370640965E-3"
$ws.Range("AK23").Value = [double]"1.9535475365449542"
$ws.Range("AL23").Value = [double]"0.41232588415750304"
$ws.Range("AM23").Value = [double]"1.0438000000000001E-2"
$ws.Range("AN23").Value = [double]"2.0458093752840187E-3"

# Row 24
$ws.Range("AG24").Value = [double]"15.686311413510827"
$ws.Range("AH24").Value = [double]"1.5974968382342152"
$ws.Range("AI24").Value = [double]"4.4472500000000005E-2"
$ws.Range("AJ24").Value = [double]"1.7122767494384731E-3"
$ws.Range("AK24").Value = [double]"1.6464349889014456"
$ws.Range("AL24").Value = [double]"0.22395941688111298"
$ws.Range("AM24").Value = [double]"9.4684999999999995E-3"
$ws.Range("AN24").Value = [double]"2.1371241579905241E-3"

# Row 44
$ws.Range("AG44").Value = [double]"8.6898954639811326"
$ws.Range("AH44").Value = [double]"0.56515879905236555"
$ws.Range("AI44").Value = [double]"4.1472500000000002E-2"
$ws.Range("AJ44").Value = [double]"1.5617378141032374E-3"
$ws.Range("AK44").Value = [double]"2.2652248834184032"
$ws.Range("AL44").Value = [double]"0.33285168385544511"
$ws.Range("AM44").Value = [double]"1.0208750000000001E-2"
$ws.Range("AN44").Value = [double]"1.0099681100576056E-3"

# Row 46
$ws.Range("AI46").Value = [double]"4.2617500000000003E-2"
$ws.Range("AJ46").Value = [double]"2.4842889660155834E-3"
$ws.Range("AK46").Value = [double]"3.2536575814417046"
$ws.Range("AL46").Value = [double]"1.0270030910943657"
$ws.Range("AM46").Value = [double]"1.57425E-2"
$ws.Range("AN46").Value = [double]"5.9179522077601543E-3"

# Row 47
$ws.Range("AG47").Value = [double]"16.69503292505626"
$ws.Range("AH47").Value = [double]"1.084706323245489"
$ws.Range("AI47").Value = [double]"4.3560000000000001E-2"
$ws.Range("AJ47").Value = [double]"1.0611628841354704E-3"
$ws.Range("AK47").Value = [double]"1.7844693537727041"
$ws.Range("AL47").Value = [double]"0.39454257254064817"
$ws.Range("AM47").Value = [double]"9.8587499999999995E-3"
$ws.Range("AN47").Value = [double]"1.6013172442294716E-3"
